$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row rename (Spanish descriptive headers -> short English codes)
$ws.Cells.Item(1,1).Value = "mx_state"
$ws.Cells.Item(1,2).Value = "mx_municipality"
$ws.Cells.Item(1,3).Value = "n_matriculas"
$ws.Cells.Item(1,4).Value = "pct_matriculas"

# Title-case Spanish connector words (de/del/la/el/los/las/y) in state/municipality names
$ws.Cells.Item(7,2).Value = "Pabellón De Arteaga"
$ws.Cells.Item(8,2).Value = "Rincón De Romos"
$ws.Cells.Item(9,2).Value = "San Francisco De Los Romo"
$ws.Cells.Item(10,2).Value = "San José De Gracia"
$ws.Cells.Item(31,2).Value = "Amatenango De La Frontera"
$ws.Cells.Item(32,2).Value = "Amatenango Del Valle"
$ws.Cells.Item(40,2).Value = "Chiapa De Corzo"
$ws.Cells.Item(61,2).Value = "Salto De Agua"
$ws.Cells.Item(62,2).Value = "San Cristóbal De Las Casas"
$ws.Cells.Item(88,2).Value = "Guadalupe Y Calvo"
$ws.Cells.Item(90,2).Value = "Hidalgo Del Parral"
$ws.Cells.Item(99,2).Value = "Valle De Zaragoza"
$ws.Cells.Item(108,2).Value = "San Juan De Sabinas"
$ws.Cells.Item(121,2).Value = "Villa De Álvarez"
$ws.Cells.Item(123,1).Value = "Ciudad De México"
$ws.Cells.Item(127,2).Value = "Cuajimalpa De Morelos"
$ws.Cells.Item(149,2).Value = "Nombre De Dios"
$ws.Cells.Item(151,2).Value = "Pánuco De Coronado"
$ws.Cells.Item(156,2).Value = "San Juan Del Río"
$ws.Cells.Item(163,1).Value = "Estado De México"
$ws.Cells.Item(163,2).Value = "Acambay De Ruíz Castañeda"
$ws.Cells.Item(166,2).Value = "Almoloya De Alquisiras"
$ws.Cells.Item(171,2).Value = "Atizapán De Zaragoza"
$ws.Cells.Item(177,2).Value = "Chapa De Mota"
$ws.Cells.Item(181,2).Value = "Coacalco De Berriozábal"
$ws.Cells.Item(185,2).Value = "Ecatepec De Morelos"
$ws.Cells.Item(187,2).Value = "Ixtapan De La Sal"
$ws.Cells.Item(188,2).Value = "Ixtapan Del Oro"
$ws.Cells.Item(198,2).Value = "Naucalpan De Juárez"
$ws.Cells.Item(206,2).Value = "San Felipe Del Progreso"
$ws.Cells.Item(208,2).Value = "San Simón De Guerrero"
$ws.Cells.Item(209,2).Value = "Soyaniquilpan De Juárez"
$ws.Cells.Item(218,2).Value = "Tenango Del Valle"
$ws.Cells.Item(228,2).Value = "Tlalnepantla De Baz"
$ws.Cells.Item(233,2).Value = "Valle De Bravo"
$ws.Cells.Item(234,2).Value = "Valle De Chalco Solidaridad"
$ws.Cells.Item(235,2).Value = "Villa De Allende"
$ws.Cells.Item(246,2).Value = "San Miguel De Allende"
$ws.Cells.Item(247,2).Value = "Apaseo El Grande"
$ws.Cells.Item(252,2).Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Cells.Item(262,2).Value = "Purísima Del Rincón"
$ws.Cells.Item(266,2).Value = "San Diego De La Unión"
$ws.Cells.Item(268,2).Value = "San Francisco Del Rincón"
$ws.Cells.Item(270,2).Value = "San Luis De La Paz"
$ws.Cells.Item(271,2).Value = "Santa Cruz De Juventino Rosas"
$ws.Cells.Item(272,2).Value = "Silao De La Victoria"
$ws.Cells.Item(277,2).Value = "Valle De Santiago"
$ws.Cells.Item(283,2).Value = "Acapulco De Juárez"
$ws.Cells.Item(285,2).Value = "Ajuchitlán Del Progreso"
$ws.Cells.Item(286,2).Value = "Alcozauca De Guerrero"
$ws.Cells.Item(289,2).Value = "Atenango Del Río"
$ws.Cells.Item(290,2).Value = "Atlamajalcingo Del Monte"
$ws.Cells.Item(292,2).Value = "Atoyac De Álvarez"
$ws.Cells.Item(293,2).Value = "Ayutla De Los Libres"
$ws.Cells.Item(296,2).Value = "Buenavista De Cuéllar"
$ws.Cells.Item(297,2).Value = "Chilapa De Álvarez"
$ws.Cells.Item(298,2).Value = "Chilpancingo De Los Bravo"
$ws.Cells.Item(299,2).Value = "Coahuayutla De José María Izazaga"
$ws.Cells.Item(304,2).Value = "Coyuca De Benítez"
$ws.Cells.Item(305,2).Value = "Coyuca De Catalán"
$ws.Cells.Item(309,2).Value = "Cuetzala Del Progreso"
$ws.Cells.Item(310,2).Value = "Cutzamala De Pinzón"
$ws.Cells.Item(316,2).Value = "Huitzuco De Los Figueroa"
$ws.Cells.Item(317,2).Value = "Iguala De La Independencia"
$ws.Cells.Item(319,2).Value = "Ixcateopan De Cuauhtémoc"
$ws.Cells.Item(320,2).Value = "Zihuatanejo De Azueta"
$ws.Cells.Item(322,2).Value = "La Unión De Isidoro Montes De Oca"
$ws.Cells.Item(325,2).Value = "Mártir De Cuilapan"
$ws.Cells.Item(338,2).Value = "Taxco De Alarcón"
$ws.Cells.Item(340,2).Value = "Técpan De Galeana"
$ws.Cells.Item(342,2).Value = "Tepecoacuilco De Trujano"
$ws.Cells.Item(344,2).Value = "Tixtla De Guerrero"
$ws.Cells.Item(348,2).Value = "Tlalixtaquilla De Maldonado"
$ws.Cells.Item(349,2).Value = "Tlapa De Comonfort"
$ws.Cells.Item(363,2).Value = "Atotonilco De Tula"
$ws.Cells.Item(364,2).Value = "Atotonilco El Grande"
$ws.Cells.Item(370,2).Value = "Cuautepec De Hinojosa"
$ws.Cells.Item(374,2).Value = "Huasca De Ocampo"
$ws.Cells.Item(377,2).Value = "Huejutla De Reyes"
$ws.Cells.Item(380,2).Value = "Jacala De Ledezma"
$ws.Cells.Item(385,2).Value = "Mineral De La Reforma"
$ws.Cells.Item(386,2).Value = "Mineral Del Chico"
$ws.Cells.Item(387,2).Value = "Mineral Del Monte"
$ws.Cells.Item(388,2).Value = "Mixquiahuala De Juárez"
$ws.Cells.Item(390,2).Value = "Omitlán De Juárez"
$ws.Cells.Item(391,2).Value = "Pachuca De Soto"
$ws.Cells.Item(392,2).Value = "Progreso De Obregón"
$ws.Cells.Item(398,2).Value = "Santiago De Anaya"
$ws.Cells.Item(402,2).Value = "Tenango De Doria"
$ws.Cells.Item(404,2).Value = "Tepehuacán De Guerrero"
$ws.Cells.Item(407,2).Value = "Tezontepec De Aldama"
$ws.Cells.Item(412,2).Value = "Tula De Allende"
$ws.Cells.Item(413,2).Value = "Tulancingo De Bravo"
$ws.Cells.Item(414,2).Value = "Zacualtipán De Ángeles"
$ws.Cells.Item(419,2).Value = "Acatlán De Juárez"
$ws.Cells.Item(420,2).Value = "Ahualulco De Mercado"
$ws.Cells.Item(424,2).Value = "Atemajac De Brizuela"
$ws.Cells.Item(427,2).Value = "Atotonilco El Alto"
$ws.Cells.Item(429,2).Value = "Autlán De Navarro"
$ws.Cells.Item(441,2).Value = "Concepción De Buenos Aires"
$ws.Cells.Item(442,2).Value = "Cuautitlán De García Barragán"
$ws.Cells.Item(449,2).Value = "Encarnación De Díaz"
$ws.Cells.Item(456,2).Value = "Huejuquilla El Alto"
$ws.Cells.Item(457,2).Value = "Ixtlahuacán De Los Membrillos"
$ws.Cells.Item(458,2).Value = "Ixtlahuacán Del Río"
$ws.Cells.Item(461,2).Value = "Jilotlán De Los Dolores"
$ws.Cells.Item(467,2).Value = "La Manzanilla De La Paz"
$ws.Cells.Item(468,2).Value = "Lagos De Moreno"
$ws.Cells.Item(475,2).Value = "Ojuelos De Jalisco"
$ws.Cells.Item(480,2).Value = "San Cristóbal De La Barranca"
$ws.Cells.Item(481,2).Value = "San Diego De Alejandría"
$ws.Cells.Item(483,2).Value = "San Juan De Los Lagos"
$ws.Cells.Item(486,2).Value = "San Martín De Bolaños"
$ws.Cells.Item(488,2).Value = "San Miguel El Alto"
$ws.Cells.Item(489,2).Value = "San Sebastián Del Oeste"
$ws.Cells.Item(490,2).Value = "Santa María De Los Ángeles"
$ws.Cells.Item(493,2).Value = "Talpa De Allende"
$ws.Cells.Item(494,2).Value = "Tamazula De Gordiano"
$ws.Cells.Item(500,2).Value = "Teocuitatlán De Corona"
$ws.Cells.Item(501,2).Value = "Tepatitlán De Morelos"
$ws.Cells.Item(504,2).Value = "Tizapán El Alto"
$ws.Cells.Item(505,2).Value = "Tlajomulco De Zúñiga"
$ws.Cells.Item(515,2).Value = "Unión De San Antonio"
$ws.Cells.Item(516,2).Value = "Unión De Tula"
$ws.Cells.Item(517,2).Value = "Valle De Guadalupe"
$ws.Cells.Item(522,2).Value = "Zacoalco De Torres"
$ws.Cells.Item(525,2).Value = "Zapotitlán De Vadillo"
$ws.Cells.Item(526,2).Value = "Zapotlán El Grande"
$ws.Cells.Item(550,2).Value = "Coalcomán De Vázquez Pallares"
$ws.Cells.Item(552,2).Value = "Cojumatlán De Régules"
$ws.Cells.Item(617,2).Value = "Tiquicheo De Nicolás Romero"
$ws.Cells.Item(642,2).Value = "Coatlán Del Río"
$ws.Cells.Item(651,2).Value = "Puente De Ixtla"
$ws.Cells.Item(654,2).Value = "Tetela Del Volcán"
$ws.Cells.Item(655,2).Value = "Tlaltizapán De Zapata"
$ws.Cells.Item(665,2).Value = "Amatlán De Cañas"
$ws.Cells.Item(666,2).Value = "Bahía De Banderas"
$ws.Cells.Item(670,2).Value = "Ixtlán Del Río"
$ws.Cells.Item(677,2).Value = "Santa María Del Oro"
$ws.Cells.Item(689,2).Value = "Lampazos De Naranjo"
$ws.Cells.Item(693,2).Value = "San Nicolás De Los Garza"
$ws.Cells.Item(696,2).Value = "Acatlán De Pérez Figueroa"
$ws.Cells.Item(700,2).Value = "Ayoquezco De Aldama"
$ws.Cells.Item(703,2).Value = "Chalcatongo De Hidalgo"
$ws.Cells.Item(705,2).Value = "Coicoyán De Las Flores"
$ws.Cells.Item(706,2).Value = "Constancia Del Rosario"
$ws.Cells.Item(708,2).Value = "Cuilápam De Guerrero"
$ws.Cells.Item(709,2).Value = "Guadalupe De Ramírez"
$ws.Cells.Item(710,2).Value = "Guevea De Humboldt"
$ws.Cells.Item(711,2).Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Cells.Item(712,2).Value = "Heroica Ciudad De Huajuapan De León"
$ws.Cells.Item(713,2).Value = "Heroica Ciudad De Tlaxiaco"
$ws.Cells.Item(715,2).Value = "Ixtlán De Juárez"
$ws.Cells.Item(716,2).Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Cells.Item(720,2).Value = "Mariscala De Juárez"
$ws.Cells.Item(721,2).Value = "Mártires De Tacubaya"
$ws.Cells.Item(724,2).Value = "Miahuatlán De Porfirio Díaz"
$ws.Cells.Item(726,2).Value = "Nejapa De Madero"
$ws.Cells.Item(727,2).Value = "Oaxaca De Juárez"
$ws.Cells.Item(728,2).Value = "Ocotlán De Morelos"
$ws.Cells.Item(729,2).Value = "Pinotepa De Don Luis"
$ws.Cells.Item(730,2).Value = "Putla Villa De Guerrero"
$ws.Cells.Item(732,2).Value = "San Agustín De Las Juntas"
$ws.Cells.Item(749,2).Value = "San Francisco Del Mar"
$ws.Cells.Item(756,2).Value = "San José Del Peñasco"
$ws.Cells.Item(761,2).Value = "San Juan Bautista Lo De Soto"
$ws.Cells.Item(766,2).Value = "San Juan Del Río"
$ws.Cells.Item(786,2).Value = "San Mateo Del Mar"
$ws.Cells.Item(791,2).Value = "San Miguel Del Puerto"
$ws.Cells.Item(816,2).Value = "Santa Ana Del Valle"
$ws.Cells.Item(822,2).Value = "Santa Cruz De Bravo"
$ws.Cells.Item(825,2).Value = "Santa Cruz Tacache De Mina"
$ws.Cells.Item(830,2).Value = "Santa Inés De Zaragoza"
$ws.Cells.Item(831,2).Value = "Santa Inés Del Monte"
$ws.Cells.Item(839,2).Value = "Santa María Jalapa Del Marqués"
$ws.Cells.Item(882,2).Value = "Tamazulápam Del Espíritu Santo"
$ws.Cells.Item(883,2).Value = "Tataltepec De Valdés"
$ws.Cells.Item(884,2).Value = "Teococuilco De Marcos Pérez"
$ws.Cells.Item(885,2).Value = "Teotitlán De Flores Magón"
$ws.Cells.Item(886,2).Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Cells.Item(887,2).Value = "Tlacolula De Matamoros"
$ws.Cells.Item(888,2).Value = "Totontepec Villa De Morelos"
$ws.Cells.Item(889,2).Value = "Villa De Chilapa De Díaz"
$ws.Cells.Item(890,2).Value = "Villa De Etla"
$ws.Cells.Item(891,2).Value = "Villa De Tututepec"
$ws.Cells.Item(892,2).Value = "Villa De Zaachila"
$ws.Cells.Item(893,2).Value = "Villa Sola De Vega"
$ws.Cells.Item(894,2).Value = "Villa Talea De Castro"
$ws.Cells.Item(897,2).Value = "Zapotitlán Del Río"
$ws.Cells.Item(899,2).Value = "Zimatlán De Álvarez"
$ws.Cells.Item(914,2).Value = "Ayotoxco De Guerrero"
$ws.Cells.Item(916,2).Value = "Chalchicomula De Sesma"
$ws.Cells.Item(930,2).Value = "Cuayuca De Andrade"
$ws.Cells.Item(931,2).Value = "Cuetzalan Del Progreso"
$ws.Cells.Item(939,2).Value = "Huehuetlán El Grande"
$ws.Cells.Item(946,2).Value = "Izúcar De Matamoros"
$ws.Cells.Item(952,2).Value = "Los Reyes De Juárez"
$ws.Cells.Item(962,2).Value = "Palmar De Bravo"
$ws.Cells.Item(976,2).Value = "San Nicolás De Los Ranchos"
$ws.Cells.Item(979,2).Value = "San Salvador El Verde"
$ws.Cells.Item(990,2).Value = "Tepatlaxco De Hidalgo"
$ws.Cells.Item(994,2).Value = "Tepexi De Rodríguez"
$ws.Cells.Item(996,2).Value = "Tepeyahualco De Cuauhtémoc"
$ws.Cells.Item(997,2).Value = "Tetela De Ocampo"
$ws.Cells.Item(1001,2).Value = "Tlacotepec De Benito Juárez"
$ws.Cells.Item(1029,2).Value = "Amealco De Bonfil"
$ws.Cells.Item(1030,2).Value = "Cadereyta De Montes"
$ws.Cells.Item(1035,2).Value = "Jalpan De Serra"
$ws.Cells.Item(1036,2).Value = "Landa De Matamoros"
$ws.Cells.Item(1039,2).Value = "Pinal De Amoles"
$ws.Cells.Item(1042,2).Value = "San Juan Del Río"
$ws.Cells.Item(1049,2).Value = "Armadillo De Los Infante"
$ws.Cells.Item(1052,2).Value = "Ciudad Del Maíz"
$ws.Cells.Item(1058,2).Value = "Mexquitic De Carmona"
$ws.Cells.Item(1062,2).Value = "San Ciro De Acosta"
$ws.Cells.Item(1065,2).Value = "Santa María Del Río"
$ws.Cells.Item(1067,2).Value = "Soledad De Graciano Sánchez"
$ws.Cells.Item(1071,2).Value = "Villa De Ramos"
$ws.Cells.Item(1103,2).Value = "Nacozari De García"
$ws.Cells.Item(1117,2).Value = "Jalpa De Méndez"
$ws.Cells.Item(1135,2).Value = "Soto La Marina"
$ws.Cells.Item(1149,2).Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Cells.Item(1150,2).Value = "Nanacamilpa De Mariano Arista"
$ws.Cells.Item(1156,2).Value = "Tetla De La Solidaridad"
$ws.Cells.Item(1167,2).Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Cells.Item(1169,2).Value = "Amatlán De Los Reyes"
$ws.Cells.Item(1175,2).Value = "Boca Del Río"
$ws.Cells.Item(1189,2).Value = "Cosamaloapan De Carpio"
$ws.Cells.Item(1190,2).Value = "Cosautlán De Carvajal"
$ws.Cells.Item(1205,2).Value = "Hueyapan De Ocampo"
$ws.Cells.Item(1206,2).Value = "Ignacio De La Llave"
$ws.Cells.Item(1209,2).Value = "Ixhuatlán De Madero"
$ws.Cells.Item(1210,2).Value = "Ixhuatlán Del Sureste"
$ws.Cells.Item(1218,2).Value = "Juchique De Ferrer"
$ws.Cells.Item(1223,2).Value = "Lerdo De Tejada"
$ws.Cells.Item(1226,2).Value = "Martínez De La Torre"
$ws.Cells.Item(1227,2).Value = "Medellín De Bravo"
$ws.Cells.Item(1230,2).Value = "Mixtla De Altamirano"
$ws.Cells.Item(1232,2).Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Cells.Item(1237,2).Value = "Ozuluama De Mascareñas"
$ws.Cells.Item(1241,2).Value = "Paso Del Macho"
$ws.Cells.Item(1244,2).Value = "Poza Rica De Hidalgo"
$ws.Cells.Item(1251,2).Value = "Sayula De Alemán"
$ws.Cells.Item(1253,2).Value = "Soledad De Doblado"
$ws.Cells.Item(1279,2).Value = "Vega De Alatorre"
$ws.Cells.Item(1314,2).Value = "Cañitas De Felipe Pescador"
$ws.Cells.Item(1323,2).Value = "Jiménez Del Teul"
$ws.Cells.Item(1326,2).Value = "Mezquital Del Oro"
$ws.Cells.Item(1330,2).Value = "Moyahua De Estrada"
$ws.Cells.Item(1331,2).Value = "Nochistlán De Mejía"
$ws.Cells.Item(1332,2).Value = "Noria De Ángeles"
$ws.Cells.Item(1340,2).Value = "Teúl De González Ortega"
$ws.Cells.Item(1341,2).Value = "Tlaltenango De Sánchez Román"
$ws.Cells.Item(1342,2).Value = "Trinidad García De La Cadena"
$ws.Cells.Item(1344,2).Value = "Villa De Cos"

# Floating point re-serialization deltas (1 ULP) from recalculated percentages
$ws.Cells.Item(120,4).Value = 0.009235936188077248
$ws.Cells.Item(147,4).Value = 0.00093292284728053
$ws.Cells.Item(172,4).Value = 0.00093292284728053
$ws.Cells.Item(187,4).Value = 0.00093292284728053
$ws.Cells.Item(214,4).Value = 0.00093292284728053
$ws.Cells.Item(291,4).Value = 0.00093292284728053
$ws.Cells.Item(306,4).Value = 0.00093292284728053
$ws.Cells.Item(378,4).Value = 0.00093292284728053
$ws.Cells.Item(407,4).Value = 0.00093292284728053
$ws.Cells.Item(444,4).Value = 0.00093292284728053
$ws.Cells.Item(469,4).Value = 0.00093292284728053
$ws.Cells.Item(492,4).Value = 0.00093292284728053
$ws.Cells.Item(496,4).Value = 0.00093292284728053
$ws.Cells.Item(525,4).Value = 0.00093292284728053
$ws.Cells.Item(528,4).Value = 0.09291911558914076
$ws.Cells.Item(529,4).Value = 0.00093292284728053
$ws.Cells.Item(651,4).Value = 0.00093292284728053
$ws.Cells.Item(681,4).Value = 0.00093292284728053
$ws.Cells.Item(802,4).Value = 0.00093292284728053
$ws.Cells.Item(860,4).Value = 0.009982274465901672
$ws.Cells.Item(938,4).Value = 0.00093292284728053
$ws.Cells.Item(1091,4).Value = 0.00093292284728053
$ws.Cells.Item(1236,4).Value = 0.00093292284728053
$ws.Cells.Item(1250,4).Value = 0.00093292284728053

# Remove trailing metadata/footnote rows (1353:1357); dimension shrinks to A1:D1351
$ws.Rows("1353:1357").Delete()